# Weekly fruit/vegetable price update: insert a new daily record for
# "Feria Lagunitas de Puerto Montt" - Lechuga (Escarola, Primera) dated
# 2021-11-11 (serial 44511), pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 305 - everything from the old row 305 onward
# shifts down by one row (old 305 -> 306, ..., old 378 -> 379).
$ws.Rows.Item(305).Insert()

# Populate the newly-inserted row 305 with the new record.
$ws.Range("A305").Value = 4
$ws.Range("B305").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C305").Value = "Los Lagos"
$ws.Range("D305").Value2 = 44511
$ws.Range("E305").Value = 10
$ws.Range("F305").Value = 100112033
$ws.Range("G305").Value = "Lechuga"
$ws.Range("H305").Value = "Escarola"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 300
$ws.Range("K305").Value = 8500
$ws.Range("L305").Value = 9000
$ws.Range("M305").Value = 8750
$ws.Range("N305").Value = "$/caja 15 unidades"
$ws.Range("O305").Value = "Región de Coquimbo"
$ws.Range("P305").Value = 583
$ws.Range("Q305").Value = 15
$ws.Range("R305").Value = "Hortaliza"
